$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: quantity/purchase_price/sale_price were entered as text; convert
# them to real numbers (matches every other data row in the sheet). The
# creation_date timestamp is also re-saved (sub-millisecond serial drift).
$ws.Range("D10").Value = 50
$ws.Range("E10").Value = 4500
$ws.Range("F10").Value = 6600
$ws.Range("G10").Value = 45813.87414978009

# Row 11: new product entry (Metrotexato). Numeric-looking columns are
# recorded as text here (consistent with how row 10 originally looked
# before its own later clean-up), so force text formatting while typing
# them in, then drop the temporary formatting so the cell keeps the
# default (unstyled) look.
$ws.Range("A11").Value = "2M"
$ws.Range("B11").Value = "Medicamentos"
$ws.Range("C11").Value = "Metrotexato 1 caja 10 pastillas 200mg"

$ws.Range("D11:F11").NumberFormat = "@"
$ws.Range("D11").Value = "20"
$ws.Range("E11").Value = "23400"
$ws.Range("F11").Value = "28000"
$ws.Range("D11:F11").ClearFormats()

$ws.Range("G11").Value = 45815.82747358619
$ws.Range("G11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
